$wb = $excel.ActiveWorkbook

# The "int" / "string" type-row entries on the HERO and MONSTER unit
# tables change from "int" to "string" for the unit_id column (A3).
$wsHero = $wb.Worksheets.Item("UNIT_HERO")
$wsMonster = $wb.Worksheets.Item("UNIT_MONSTER")

$wsHero.Range("A3").Value = "string"
$wsMonster.Range("A3").Value = "string"

# Switch the active tab from UNIT_LEVEL@HERO back to UNIT_HERO, and move
# the selection on UNIT_HERO from A5 to A3.
$wsHero.Activate()
$null = $wsHero.Range("A3").Select()
